# Scheduled-runner refresh of market/profit figures on the Leve-profit
# sheets (currentAveragePrice* / LevePrice* / LeveProfit* columns).
# Updates 37 rows across all 8 job sheets with refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1196.7778
$ws.Range("I18").Value = 471.375
$ws.Range("K18").Value = 471.375
$ws.Range("M18").Value = -187.375

$ws.Range("H64").Value = 2407639.8
$ws.Range("I64").Value = 3849847
$ws.Range("J64").Value = 3961.0833
$ws.Range("K64").Value = 3849847
$ws.Range("L64").Value = 3961.0833
$ws.Range("M64").Value = -3849599
$ws.Range("N64").Value = -4457.0833

$ws.Range("H67").Value = 2407639.8
$ws.Range("I67").Value = 3849847
$ws.Range("J67").Value = 3961.0833
$ws.Range("K67").Value = 3849847
$ws.Range("L67").Value = 3961.0833
$ws.Range("M67").Value = -3848989
$ws.Range("N67").Value = -5677.0833

$ws.Range("H76").Value = 4102.4375
$ws.Range("I76").Value = 4012.6365
$ws.Range("J76").Value = 4300
$ws.Range("K76").Value = 4012.6365
$ws.Range("L76").Value = 4300
$ws.Range("M76").Value = -3697.6365
$ws.Range("N76").Value = -4930

$ws.Range("H79").Value = 4102.4375
$ws.Range("I79").Value = 4012.6365
$ws.Range("J79").Value = 4300
$ws.Range("K79").Value = 4012.6365
$ws.Range("L79").Value = 4300
$ws.Range("M79").Value = -2920.6365
$ws.Range("N79").Value = -6484

$ws.Range("H98").Value = 1166.4286
$ws.Range("I98").Value = 860.8333
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 860.8333
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 637.1667
$ws.Range("N98").Value = -5996

$ws.Range("H122").Value = 1166.4286
$ws.Range("I122").Value = 860.8333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2582.4999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -132.4998999999998
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 5558546
$ws.Range("I132").Value = 2430.2432
$ws.Range("J132").Value = 31255580
$ws.Range("K132").Value = 7290.7296
$ws.Range("L132").Value = 93766740
$ws.Range("M132").Value = -4760.7296
$ws.Range("N132").Value = -93771800

$ws.Range("H137").Value = 2581.4
$ws.Range("I137").Value = 3001
$ws.Range("J137").Value = 903
$ws.Range("K137").Value = 9003
$ws.Range("L137").Value = 2709
$ws.Range("M137").Value = -6453
$ws.Range("N137").Value = -7809

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5603.87
$ws.Range("I32").Value = 4360.6406
$ws.Range("J32").Value = 15662.728
$ws.Range("K32").Value = 4360.6406
$ws.Range("L32").Value = 15662.728
$ws.Range("M32").Value = -4073.6406
$ws.Range("N32").Value = -16236.728

$ws.Range("H63").Value = 125004850
$ws.Range("I63").Value = 166671140
$ws.Range("J63").Value = 6006
$ws.Range("K63").Value = 166671140
$ws.Range("L63").Value = 6006
$ws.Range("M63").Value = -166670454
$ws.Range("N63").Value = -7378

$ws.Range("H66").Value = 125004850
$ws.Range("I66").Value = 166671140
$ws.Range("J66").Value = 6006
$ws.Range("K66").Value = 833355700
$ws.Range("L66").Value = 30030
$ws.Range("M66").Value = -833352268
$ws.Range("N66").Value = -36894

$ws.Range("H88").Value = 1885.6666
$ws.Range("J88").Value = 1828.5
$ws.Range("L88").Value = 1828.5
$ws.Range("N88").Value = -2640.5

$ws.Range("H91").Value = 1885.6666
$ws.Range("J91").Value = 1828.5
$ws.Range("L91").Value = 1828.5
$ws.Range("N91").Value = -4636.5

$ws.Range("H132").Value = 6580782.5
$ws.Range("I132").Value = 8622467
$ws.Range("J132").Value = 2019.1111
$ws.Range("K132").Value = 25867401
$ws.Range("L132").Value = 6057.3333
$ws.Range("M132").Value = -25864871
$ws.Range("N132").Value = -11117.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4047.8484
$ws.Range("I105").Value = 2389.875
$ws.Range("J105").Value = 4578.4
$ws.Range("K105").Value = 2389.875
$ws.Range("L105").Value = 4578.4
$ws.Range("M105").Value = -642.875
$ws.Range("N105").Value = -8072.4

$ws.Range("H118").Value = 7890
$ws.Range("J118").Value = 7890
$ws.Range("L118").Value = 7890
$ws.Range("N118").Value = -11204

$ws.Range("H134").Value = 3168.077
$ws.Range("I134").Value = 1732.2903
$ws.Range("J134").Value = 5287.5713
$ws.Range("K134").Value = 5196.8709
$ws.Range("L134").Value = 15862.7139
$ws.Range("M134").Value = -2661.8709
$ws.Range("N134").Value = -20932.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2017.7975
$ws.Range("I31").Value = 1724.1957
$ws.Range("J31").Value = 2427.0605
$ws.Range("K31").Value = 1724.1957
$ws.Range("L31").Value = 2427.0605
$ws.Range("M31").Value = -1429.1957
$ws.Range("N31").Value = -3017.0605

$ws.Range("H34").Value = 2017.7975
$ws.Range("I34").Value = 1724.1957
$ws.Range("J34").Value = 2427.0605
$ws.Range("K34").Value = 1724.1957
$ws.Range("L34").Value = 2427.0605
$ws.Range("M34").Value = -1522.1957
$ws.Range("N34").Value = -2831.0605

$ws.Range("H62").Value = 2520
$ws.Range("I62").Value = 2340
$ws.Range("J62").Value = 2571.4285
$ws.Range("K62").Value = 2340
$ws.Range("L62").Value = 2571.4285
$ws.Range("M62").Value = -1716
$ws.Range("N62").Value = -3819.4285

$ws.Range("H65").Value = 2520
$ws.Range("I65").Value = 2340
$ws.Range("J65").Value = 2571.4285
$ws.Range("K65").Value = 11700
$ws.Range("L65").Value = 12857.1425
$ws.Range("M65").Value = -8580
$ws.Range("N65").Value = -19097.1425

$ws.Range("H134").Value = 2839.258
$ws.Range("I134").Value = 2691.7273
$ws.Range("J134").Value = 3199.889
$ws.Range("K134").Value = 8075.1819
$ws.Range("L134").Value = 9599.667000000001
$ws.Range("M134").Value = -5540.1819
$ws.Range("N134").Value = -14669.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3873.9546
$ws.Range("I109").Value = 1854.5
$ws.Range("J109").Value = 4631.25
$ws.Range("K109").Value = 5563.5
$ws.Range("L109").Value = 13893.75
$ws.Range("M109").Value = -4523.5
$ws.Range("N109").Value = -15973.75

$ws.Range("H113").Value = 45455320
$ws.Range("I113").Value = 142857700
$ws.Range("J113").Value = 876.6667
$ws.Range("K113").Value = 428573100
$ws.Range("L113").Value = 2630.0001
$ws.Range("M113").Value = -428570930
$ws.Range("N113").Value = -6970.0001

$ws.Range("H131").Value = 1194.25
$ws.Range("J131").Value = 1275.4878
$ws.Range("L131").Value = 3826.463400000001
$ws.Range("N131").Value = -13906.4634

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13048.591
$ws.Range("I70").Value = 23767.1
$ws.Range("J70").Value = 4116.5
$ws.Range("K70").Value = 23767.1
$ws.Range("L70").Value = 4116.5
$ws.Range("M70").Value = -23497.1
$ws.Range("N70").Value = -4656.5

$ws.Range("H73").Value = 13048.591
$ws.Range("I73").Value = 23767.1
$ws.Range("J73").Value = 4116.5
$ws.Range("K73").Value = 23767.1
$ws.Range("L73").Value = 4116.5
$ws.Range("M73").Value = -22831.1
$ws.Range("N73").Value = -5988.5

$ws.Range("H80").Value = 22225176
$ws.Range("I80").Value = 83335410
$ws.Range("J80").Value = 3272.7273
$ws.Range("K80").Value = 83335410
$ws.Range("L80").Value = 3272.7273
$ws.Range("M80").Value = -83334412
$ws.Range("N80").Value = -5268.7273

$ws.Range("H83").Value = 22225176
$ws.Range("I83").Value = 83335410
$ws.Range("J83").Value = 3272.7273
$ws.Range("K83").Value = 416677050
$ws.Range("L83").Value = 16363.6365
$ws.Range("M83").Value = -416672058
$ws.Range("N83").Value = -26347.6365

$ws.Range("H126").Value = 4279.6206
$ws.Range("I126").Value = 3351.2
$ws.Range("J126").Value = 4768.263
$ws.Range("K126").Value = 10053.6
$ws.Range("L126").Value = 14304.789
$ws.Range("M126").Value = -7583.599999999999
$ws.Range("N126").Value = -19244.789

$ws.Range("H132").Value = 4728.1714
$ws.Range("I132").Value = 4805.6772
$ws.Range("K132").Value = 14417.0316
$ws.Range("M132").Value = -11887.0316

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6163
$ws.Range("I40").Value = 8913.571
$ws.Range("J40").Value = 4558.5
$ws.Range("K40").Value = 8913.571
$ws.Range("L40").Value = 4558.5
$ws.Range("M40").Value = -8777.571
$ws.Range("N40").Value = -4830.5

$ws.Range("H122").Value = 7539.231
$ws.Range("I122").Value = 8870.25
$ws.Range("J122").Value = 6947.6665
$ws.Range("K122").Value = 26610.75
$ws.Range("L122").Value = 20842.9995
$ws.Range("M122").Value = -24160.75
$ws.Range("N122").Value = -25742.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5254.5454
$ws.Range("I62").Value = 5685.7144
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 5685.7144
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -5061.7144
$ws.Range("N62").Value = -5748

$ws.Range("H65").Value = 5254.5454
$ws.Range("I65").Value = 5685.7144
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 28428.572
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -25308.572
$ws.Range("N65").Value = -28740

$ws.Range("H124").Value = 30429
$ws.Range("J124").Value = 30429
$ws.Range("L124").Value = 30429
$ws.Range("N124").Value = -40249

$ws.Range("H132").Value = 1763.079
$ws.Range("I132").Value = 1712.8667
$ws.Range("J132").Value = 1951.375
$ws.Range("K132").Value = 5138.6001
$ws.Range("L132").Value = 5854.125
$ws.Range("M132").Value = -2608.6001
$ws.Range("N132").Value = -10914.125
